# Reorganize the "English" BodyParts list:
#  - Move the cardio-respiratory group (Bronch_cartilage, Bronchi, Diaphragme,
#    Heart, Lungs, Trachea) from the end of the list up to just before the
#    digestive-system section.
#  - Add a new "dog body background" group (Body, Eyes, Teeth, Tongue) right
#    after the cardio-respiratory group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Cut the existing cardio-respiratory rows (currently at rows 160-165,
#    just before the end of the sheet) and re-insert them at row 134, ahead
#    of the digestive-system block.
$cardioRespiratory = @(
    "Bronch_cartilage",
    "Bronchi",
    "Diaphragme",
    "Heart",
    "Lungs",
    "Trachea"
)

$ws.Range("A160:A165").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null

$ws.Rows.Item(134).Resize(6).Insert() | Out-Null
for ($i = 0; $i -lt $cardioRespiratory.Length; $i++) {
    $ws.Range("A" + (134 + $i)).Value = $cardioRespiratory[$i]
}

# 2) Insert the new "dog body background" rows right after the group that was
#    just relocated (new rows 140-143).
$bodyBackground = @(
    "Body",
    "Eyes",
    "Teeth",
    "Tongue"
)

$ws.Rows.Item(140).Resize(4).Insert() | Out-Null
for ($i = 0; $i -lt $bodyBackground.Length; $i++) {
    $ws.Range("A" + (140 + $i)).Value = $bodyBackground[$i]
}
